$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = 74.74473542407863
$ws.Range("D3").Value  = 74.74473542407854
$ws.Range("D4").Value  = 0.9993073109016987
$ws.Range("D5").Value  = 71.2471598359226
$ws.Range("D6").Value  = 0.9993073109016987
$ws.Range("D8").Value  = 699.515117631189
$ws.Range("D10").Value = 266.2820652073141
$ws.Range("D11").Value = 266.2820652073141
$ws.Range("D12").Value = 240.280632614111
$ws.Range("D13").Value = 240.280632614111
$ws.Range("D14").Value = 3.497575588155945
